$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-40 down to 10-41
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new data record
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44804
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 100112042
$ws.Range("G9").Value = "Locoto"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 3300
$ws.Range("L9").Value = 3300
$ws.Range("M9").Value = 3300
$ws.Range("N9").Value = "$/kilo"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 3300
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
